# Update "Orders" sheet (sheet1): extend the order list with rows 70-78,
# and update row 70's PackageID from 1 to 11.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# The new data occupies 8 brand-new rows after the existing last row (70),
# so insert them first (row 70 itself is edited in place below).
$ws.Range("A71:A78").EntireRow.Insert()

# These "numeric" values are stored as text in the workbook (matches the
# numberStoredAsText ignoredError convention already used throughout the
# sheet), so force a text number format before assigning them - otherwise
# Excel would coerce them to real numbers (dropping e.g. leading zeros).
$ws.Range("A70").NumberFormat = "@"
$ws.Range("F70:F78").NumberFormat = "@"
$ws.Range("A78").NumberFormat = "@"

$ws.Range("A70").Value = "11"
$ws.Range("C70").Value = "46_拉丝橙_Spider orange_Gerbera L._20stems"
$ws.Range("F70").Value = "10"

$ws.Range("C71").Value = "77_珍爱mini_undefined_Gerbera L._20stems"
$ws.Range("F71").Value = "20"

$ws.Range("C72").Value = "83_布拉格_undefined_Gerbera L._10stems"
$ws.Range("F72").Value = "30"

$ws.Range("C73").Value = "411_紫罗兰白_violet white_undefined_1bunch"
$ws.Range("F73").Value = "20"

$ws.Range("C74").Value = "509_翠珠粉_Didiscus caeruleus`npink_Trachymene Coerulea_1bunch"
$ws.Range("F74").Value = "10"

$ws.Range("C75").Value = "578_腊梅粉_wax pink_undefined_1bunch"
$ws.Range("F75").Value = "17"

$ws.Range("C76").Value = "354_桔叶_undefined_undefined_1bunch"
$ws.Range("F76").Value = "15"

$ws.Range("C77").Value = "550_蕾丝红色_lace flower brown_undefined_1bunch"
$ws.Range("F77").Value = "5"

$ws.Range("A78").Value = "12"
$ws.Range("C78").Value = "604_康乃馨粉佳人_pink_undefined_20stems"
$ws.Range("F78").Value = "30"

# Update "Summary" sheet (sheet2): G2 total-number string reflects the
# newly added rows. It's a 129-digit run of digits that must stay text (a
# plain numeric assignment would round-trip through a double and lose
# precision), so force it with a text number format too.
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("G2").NumberFormat = "@"
$ws2.Range("G2").Value = "020102010555553010515103010301515151512101015551018915102057128811910413511553020102053101020201551051055510151110203020101715530"
